$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H did not exist before (sheet was A1:G12); copy the header-cell
# formatting (bold font, border, center/top alignment) from G1 into H1 before
# writing values, so the new column matches the existing header style.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

# --- Header row (row 1): insert "ABG-only" as a new column between
# "ABG+VBG+OTHER" and "OTHER-only"; OTHER-only/VBG+OTHER/VBG-only shift right ---
$ws.Cells.Item(1, 1).Value = 'symptom_group'
$ws.Cells.Item(1, 2).Value = 'ABG+OTHER'
$ws.Cells.Item(1, 3).Value = 'ABG+VBG'
$ws.Cells.Item(1, 4).Value = 'ABG+VBG+OTHER'
$ws.Cells.Item(1, 5).Value = 'ABG-only'
$ws.Cells.Item(1, 6).Value = 'OTHER-only'
$ws.Cells.Item(1, 7).Value = 'VBG+OTHER'
$ws.Cells.Item(1, 8).Value = 'VBG-only'

# --- Data rows 2-12: row labels/values updated per refreshed cohort export ---
# Row 2: Diseases (patient-stated)
$ws.Cells.Item(2, 1).Value = 'Diseases (patient-stated)'
$ws.Cells.Item(2, 2).Value = 5.6
$ws.Cells.Item(2, 3).Value = 4.4
$ws.Cells.Item(2, 4).Value = 4
$ws.Cells.Item(2, 5).Value = 7.7
$ws.Cells.Item(2, 6).Value = 6.9
$ws.Cells.Item(2, 7).Value = 3.9
$ws.Cells.Item(2, 8).Value = 5.3

# Row 3: Injuries & adverse effects
$ws.Cells.Item(3, 1).Value = 'Injuries & adverse effects'
$ws.Cells.Item(3, 2).Value = 13.5
$ws.Cells.Item(3, 3).Value = 8.5
$ws.Cells.Item(3, 4).Value = 4.9
$ws.Cells.Item(3, 5).Value = 12
$ws.Cells.Item(3, 6).Value = 6.9
$ws.Cells.Item(3, 7).Value = 3.1
$ws.Cells.Item(3, 8).Value = 6.3

# Row 4: Other
$ws.Cells.Item(4, 1).Value = 'Other'
$ws.Cells.Item(4, 2).Value = 10.9
$ws.Cells.Item(4, 3).Value = 5.7
$ws.Cells.Item(4, 4).Value = 7.8
$ws.Cells.Item(4, 5).Value = 7.1
$ws.Cells.Item(4, 6).Value = 6.2
$ws.Cells.Item(4, 7).Value = 4.7
$ws.Cells.Item(4, 8).Value = 5.6

# Row 5: Symptom – Circulatory
$ws.Cells.Item(5, 1).Value = 'Symptom – Circulatory'
$ws.Cells.Item(5, 2).Value = 9.8
$ws.Cells.Item(5, 3).Value = 7.3
$ws.Cells.Item(5, 4).Value = 8
$ws.Cells.Item(5, 5).Value = 12.3
$ws.Cells.Item(5, 6).Value = 8
$ws.Cells.Item(5, 7).Value = 9.4
$ws.Cells.Item(5, 8).Value = 9.4

# Row 6: Symptom – Digestive
$ws.Cells.Item(6, 1).Value = 'Symptom – Digestive'
$ws.Cells.Item(6, 2).Value = 14.3
$ws.Cells.Item(6, 3).Value = 10
$ws.Cells.Item(6, 4).Value = 9.6
$ws.Cells.Item(6, 5).Value = 12.7
$ws.Cells.Item(6, 6).Value = 14.9
$ws.Cells.Item(6, 7).Value = 10.2
$ws.Cells.Item(6, 8).Value = 11.7

# Row 7: Symptom – General
$ws.Cells.Item(7, 1).Value = 'Symptom – General'
$ws.Cells.Item(7, 2).Value = 4.1
$ws.Cells.Item(7, 3).Value = 4.6
$ws.Cells.Item(7, 4).Value = 4.9
$ws.Cells.Item(7, 5).Value = 4.8
$ws.Cells.Item(7, 6).Value = 6.2
$ws.Cells.Item(7, 7).Value = 7.5
$ws.Cells.Item(7, 8).Value = 5.4

# Row 8: Symptom – Genitourinary
$ws.Cells.Item(8, 1).Value = 'Symptom – Genitourinary'
$ws.Cells.Item(8, 2).Value = 4.5
$ws.Cells.Item(8, 3).Value = 5
$ws.Cells.Item(8, 4).Value = 5.2
$ws.Cells.Item(8, 5).Value = 6.3
$ws.Cells.Item(8, 6).Value = 5.4
$ws.Cells.Item(8, 7).Value = 4.7
$ws.Cells.Item(8, 8).Value = 5.5

# Row 9: Symptom – Nervous
$ws.Cells.Item(9, 1).Value = 'Symptom – Nervous'
$ws.Cells.Item(9, 2).Value = 5.6
$ws.Cells.Item(9, 3).Value = 11.8
$ws.Cells.Item(9, 4).Value = 11.6
$ws.Cells.Item(9, 5).Value = 10
$ws.Cells.Item(9, 6).Value = 13.4
$ws.Cells.Item(9, 7).Value = 14.1
$ws.Cells.Item(9, 8).Value = 13.2

# Row 10: Symptom – Respiratory
$ws.Cells.Item(10, 1).Value = 'Symptom – Respiratory'
$ws.Cells.Item(10, 2).Value = 22.6
$ws.Cells.Item(10, 3).Value = 38.5
$ws.Cells.Item(10, 4).Value = 40.3
$ws.Cells.Item(10, 5).Value = 21.3
$ws.Cells.Item(10, 6).Value = 27.9
$ws.Cells.Item(10, 7).Value = 40
$ws.Cells.Item(10, 8).Value = 34.7

# Row 11: Symptom – Skin/Hair/Nails
$ws.Cells.Item(11, 1).Value = 'Symptom – Skin/Hair/Nails'
$ws.Cells.Item(11, 2).Value = 4.5
$ws.Cells.Item(11, 3).Value = 2.2
$ws.Cells.Item(11, 4).Value = 2.4
$ws.Cells.Item(11, 5).Value = 2.8
$ws.Cells.Item(11, 6).Value = 2.9
$ws.Cells.Item(11, 7).Value = 1.2
$ws.Cells.Item(11, 8).Value = 1.9

# Row 12: Uncodable/Unknown
$ws.Cells.Item(12, 1).Value = 'Uncodable/Unknown'
$ws.Cells.Item(12, 2).Value = 4.5
$ws.Cells.Item(12, 3).Value = 2
$ws.Cells.Item(12, 4).Value = 1.4
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1.4
$ws.Cells.Item(12, 7).Value = 1.2
$ws.Cells.Item(12, 8).Value = 0.9

